# Update cryptocurrency price (D) and volume change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.999.26"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "2.216.79"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.25"
$ws.Range("E5").Value = "  -2.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.51"
$ws.Range("E7").Value = "  -1.24%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -2.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.76"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0952"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.09"
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").Value = "2.549.36"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.23"
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.839"
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("D17").Value = "2.214.05"
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("D18").Value = "41.904.65"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000107"
$ws.Range("E19").Value = "  +8.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.26"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.68"
$ws.Range("E22").Value = "  +19.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.65"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("E24").Value = "  -6.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.70"
$ws.Range("E25").Value = "  +2.07%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.95"
$ws.Range("E30").Value = "  -1.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.46"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.68"
$ws.Range("E32").Value = "  +7.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0795"
$ws.Range("E33").Value = "  -3.50%  "
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "29.25"
$ws.Range("E35").Value = "  -5.57%  "
$ws.Range("E36").Value = "  -9.93%  "
$ws.Range("E37").Value = "  -5.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0300"
$ws.Range("E38").Value = "  -4.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.54"
$ws.Range("E39").Value = "  -3.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "66.31"
$ws.Range("E40").Value = "  +6.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.12"
$ws.Range("E41").Value = "  -2.92%  "
$ws.Range("E42").Value = "  -2.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.199"
$ws.Range("E43").Value = "  -3.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.71"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.12"
$ws.Range("E45").Value = "  -2.82%  "
$ws.Range("E46").Value = "  -2.66%  "
$ws.Range("E47").Value = "  +3.09%  "
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.16"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.70"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").Value = "2.424.13"
$ws.Range("E51").Value = "  -1.54%  "
